$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.411.87'
$ws.Range('E2').Value = '  -2.44%  '
$ws.Range('D3').Value = '3.310.64'
$ws.Range('E3').Value = '  -3.03%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '558.74'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -2.94%  '
$ws.Range('E6').Value = '  -3.96%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').Value = '3.309.84'
$ws.Range('E8').Value = '  -3.05%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.468'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -2.78%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.85'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -2.33%  '
$ws.Range('E11').Value = '  -3.34%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.407'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -1.03%  '
$ws.Range('D13').Value = '3.878.74'
$ws.Range('E14').Value = '  +0.16%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '26.89'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -4.95%  '
$ws.Range('D16').Value = '3.310.36'
$ws.Range('E16').Value = '  -3.00%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000165'
$ws.Range('D17').ClearFormats()
$ws.Range('D18').Value = '60.430.53'
$ws.Range('E18').Value = '  -2.42%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.13'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -3.34%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.27'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -0.80%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '8.67'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -2.09%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '374.65'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -1.42%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '74.64'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -0.58%  '
$ws.Range('E24').Value = '  -0.09%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.538'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -4.50%  '
$ws.Range('D26').Value = '3.447.77'
$ws.Range('E26').Value = '  -3.01%  '
$ws.Range('E27').Value = '  -7.32%  '
$ws.Range('E28').Value = '  -4.16%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.998'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -0.17%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.20'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -5.21%  '
$ws.Range('E31').Value = '  -0.05%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.04'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -3.19%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '7.60'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -3.57%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '22.65'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -1.62%  '
$ws.Range('E35').Value = '  -6.47%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.14'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -5.59%  '
$ws.Range('B37').Value = 'Monero'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '166.44'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -1.59%  '
$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.53'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -3.19%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.70'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -2.20%  '
$ws.Range('D40').Value = '3.341.21'
$ws.Range('E40').Value = '  -3.06%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '26.76'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -13.22%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0731'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -5.48%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '41.96'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -1.46%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.752'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -2.85%  '
$ws.Range('E45').Value = '  -4.36%  '
$ws.Range('E46').Value = '  -4.77%  '
$ws.Range('E47').Value = '  -3.75%  '
$ws.Range('D48').Value = '2.365.54'
$ws.Range('E48').Value = '  -6.74%  '
$ws.Range('E49').Value = '  -0.02%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.42'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -6.46%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '21.29'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -4.72%  '
